# DaySale report update: a new sold item is recorded, which pushes the
# running-total row and the footer row down by one row, and the footer
# timestamp is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row at position 9. This shifts the old row 9 (running
#    total) down to row 10 and the old row 10 (footer) down to row 11,
#    carrying their values, number formats and merged cells along.
# ---------------------------------------------------------------------
$ws.Rows.Item(9).Insert()

# ---------------------------------------------------------------------
# 2. The freshly inserted row 9 needs the same per-column formatting as
#    the data rows above it (row 7 / row 8). Copying row 8's formats
#    down re-creates that look (font, fill, borders, number formats)
#    using the very same style records already used by rows 7-8.
# ---------------------------------------------------------------------
$ws.Range("A8:Q8").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Fill in the new sale line (row 9) - a new item, "شفاط ثدي الجو".
# ---------------------------------------------------------------------
$ws.Range("A9").NumberFormat = "General"
$ws.Range("A9").Value = 3

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "شفاط ثدي الجو"

$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "2:0"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "0"

$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "25.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "25.0000"

$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "1:0"

# Merge the column groups in row 9 the same way rows 7-8 are merged.
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Rows.Item(9).RowHeight = 25.5

# ---------------------------------------------------------------------
# 4. The running-total row (now row 10) reflects the new item's price.
# ---------------------------------------------------------------------
$ws.Range("P10").Value = 151
$ws.Rows.Item(10).RowHeight = 24.75

# ---------------------------------------------------------------------
# 5. The footer row (now row 11) gets the refreshed generation time.
# ---------------------------------------------------------------------
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "Monday, 18 August, 2025 9:34 AM"

Write-Host "DaySale row insert applied"
